$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 (Highs-BigM (100,100) / FEASIBLE_POINT)
$ws.Range("E2").Value = 0.00027738
$ws.Range("F2").Value = 0.01492443
$ws.Range("G2").Value = 0.0005435682127240926

# Row 3 (Highs-BigM (100,100) / OPTIMAL)
$ws.Range("E3").Value = 0.57983175
$ws.Range("F3").Value = 0.59989392
$ws.Range("G3").Value = 0.58973956
